$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.211.46"
$ws.Range("E2").Value = "  +0.78%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.600.92"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "303.42"
$ws.Range("E6").Value = "  +0.69%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3764"
$ws.Range("E7").Value = "  -0.33%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "52.11"
$ws.Range("E8").Value = "  +4.58%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3629"
$ws.Range("E9").Value = "  -0.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.268"
$ws.Range("E10").Value = "  +0.70%  "
$ws.Range("E11").Value = "  +0.18%  "
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.72"
$ws.Range("E13").Value = "  +0.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.557"
$ws.Range("E14").Value = "  -0.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.393"
$ws.Range("E15").Value = "  +0.81%  "
$ws.Range("E16").Value = "  +0.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.599.67"
$ws.Range("E17").Value = "  -0.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.17"
$ws.Range("E18").Value = "  +2.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06921"
$ws.Range("E19").Value = "  +1.60%  "
$ws.Range("E20").Value = "  -0.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.518"
$ws.Range("E21").Value = "  -0.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.005"
$ws.Range("E22").Value = "  +0.36%  "
$ws.Range("E23").Value = "  -2.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.204.97"
$ws.Range("E24").Value = "  +0.78%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.434"
$ws.Range("E25").Value = "  +2.96%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.028"
$ws.Range("E26").Value = "  +7.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.15"
$ws.Range("E27").Value = "  +0.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "149.22"
$ws.Range("E28").Value = "  -1.27%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.257"
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "135.37"
$ws.Range("E30").Value = "  +1.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.378"
$ws.Range("E31").Value = "  +6.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.693"
$ws.Range("E32").Value = "  -2.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.776.30"
$ws.Range("E33").Value = "  -0.69%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9604"
$ws.Range("E34").Value = "  -0.94%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07463"
$ws.Range("E35").Value = "  -1.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.30"
$ws.Range("E36").Value = "  -1.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02736"
$ws.Range("E37").Value = "  +1.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2523"
$ws.Range("E38").Value = "  -0.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.08784"
$ws.Range("E39").Value = "  -1.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.082"
$ws.Range("E40").Value = "  -3.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.380"
$ws.Range("E41").Value = "  +0.83%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7073"
$ws.Range("E42").Value = "  +0.77%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.39"
$ws.Range("E43").Value = "  -0.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.48"
$ws.Range("E44").Value = "  +1.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6516"
$ws.Range("E45").Value = "  -1.51%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"
$ws.Range("E46").Value = "  +0.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.309"
$ws.Range("E47").Value = "  +0.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.007"
$ws.Range("E48").Value = "  +0.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "132.01"
$ws.Range("E49").Value = "  -0.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07913"
$ws.Range("E50").Value = "  +0.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.199"
$ws.Range("E51").Value = "  -1.04%  "
